# Kanban Board update:
# - "Add Books System" moves from Not Started (A5) to Doing (B4)
# - "Update Books System" (A6) is renamed to "Edit a Book System" and takes A5's spot
# - "Website UI" (A7) moves up to A6
# - "View Book Instance Page" (A8) is renamed to "View a Book Instance Page" and moves up to A7
# - A8 becomes empty (list shrank by one)
# - New task "View all Books System" is added to Doing (B5)
# - Formulas in A3/B3/C1 recalculate automatically

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Add Books System" moves from Not Started to Doing
$ws.Range("B4").Value = "Add Books System"

# "Update Books System" is renamed/replaced by "Edit a Book System" (stays in Not Started)
$ws.Range("A5").Value = "Edit a Book System"

# "Website UI" shifts up from A7 to A6
$ws.Range("A6").Value = "Website UI"

# "View Book Instance Page" renamed to "View a Book Instance Page", shifts up to A7
$ws.Range("A7").Value = "View a Book Instance Page"

# New "Doing" task
$ws.Range("B5").Value = "View all Books System"

# A8 is now empty
$ws.Range("A8").Value = ""

# Update the active selection to B10
$ws.Range("B10").Select()
